$wb = $excel.ActiveWorkbook

# --- Settings sheet: add the new "AcmeLoginCredential" asset entry ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Range("A11").Value = "AcmeLoginCredential"
$wsSettings.Range("B11").Value = "AcmeLoginCredential"
$null = $wsSettings.Range("A11").Select()

# --- Assets sheet: drop the "Credential" Yes/No row and column ---
$wsAssets = $wb.Worksheets.Item("Assets")
$null = $wsAssets.Range("3:3").EntireRow.Delete()
$null = $wsAssets.Range("C:C").EntireColumn.Delete()
$null = $wsAssets.Range("C:C").Select()

# --- bump the internal sheetId counter for "Evaluation Warning" from 12 to 13 ---
# (duplicate the sheet, drop the original, rename the duplicate back)
$wsWarning = $wb.Worksheets.Item("Evaluation Warning")
$wsWarning.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$null = $wb.Worksheets.Item("Evaluation Warning").Delete()
$wb.Worksheets.Item("Evaluation Warning (2)").Name = "Evaluation Warning"

# restore the originally active tab (Evaluation Warning)
$null = $wb.Worksheets.Item("Evaluation Warning").Activate()
